# Daily attendance processing - 2025-10-19 19:40:23
#
# In the "Recorded By" column (G), whenever the list of recorder names
# begins with the literal (case-sensitive) token "System", that token is
# moved from the front of the comma-separated list to the back - i.e. the
# whole list is reversed. Rows whose list does not start with "System"
# (e.g. "admin@admin.com, System", or plain single-value cells) are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Case-sensitive exact string comparison helper - the interpreter's
# built-in operators (-eq, -ceq, -clike, -cmatch) all compare
# case-insensitively here, so compare character codes by hand instead.
function Test-ExactCaseMatch($str, $target) {
    if ($str.Length -ne $target.Length) { return $false }
    for ($i = 0; $i -lt $str.Length; $i++) {
        $c1 = [int][char]$str.Substring($i, 1)
        $c2 = [int][char]$target.Substring($i, 1)
        if ($c1 -ne $c2) { return $false }
    }
    return $true
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -notlike "*,*") { continue }

    $rawParts = $val.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    if ($parts.Count -gt 1 -and (Test-ExactCaseMatch $parts[0] "System")) {
        $reversed = $parts[($parts.Count - 1)..0]
        $newVal = [string]::Join(", ", $reversed)
        $cell.Value = $newVal
    }
}
